$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing hyperlinks on D2 (piyush) / D3 (shaweta) before the row
# shift -- they get re-created (in the right order, with a new row for
# "Parth Gupta" ahead of them) further down.
$ws.Range("D2").Hyperlinks.Delete()
$ws.Range("D3").Hyperlinks.Delete()

# Stash the original hyperlink-cell font style (blue, non-underlined) in a
# scratch cell so it can be re-applied after new hyperlinks are created
# (Hyperlinks.Add swaps in the builtin "Hyperlink" cell style otherwise).
$ws.Range("D2").Copy($ws.Range("Z1"))

# Insert a new row above current row 2 ("Piyush Sharma"), pushing the
# existing two data rows down to rows 3 and 4.
$ws.Rows.Item(2).Insert()

# Give the freshly-inserted row the same row height as the other data rows.
$ws.Rows.Item(2).RowHeight = 14.9

# Apply the stashed hyperlink-cell style to all three email cells.
$ws.Range("Z1").Copy($ws.Range("D2"))
$ws.Range("Z1").Copy($ws.Range("D3"))
$ws.Range("Z1").Copy($ws.Range("D4"))

# --- Row 2: new student, Parth Gupta -------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 12528
$ws.Range("C2").Value = "Parth Gupta"
$ws.Range("D2").Value = "parth@gmail.com"
$ws.Range("E2").Value = "CSE"
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 80
$ws.Range("H2").Value = 8

# --- Row 3: Piyush Sharma (renumbered sno, same data as before) ----------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 12502
$ws.Range("C3").Value = "Piyush Sharma"
$ws.Range("D3").Value = "piyush@gmail.com"
$ws.Range("E3").Value = "CSE"
$ws.Range("F3").Value = 95
$ws.Range("G3").Value = 95
$ws.Range("H3").Value = 9

# --- Row 4: Shaweta Choudhary (renumbered sno, same data as before) ------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 12518
$ws.Range("C4").Value = "Shaweta Choudhary"
$ws.Range("D4").Value = "shaweta@gmail.com"
$ws.Range("E4").Value = "CSE"
$ws.Range("F4").Value = 90
$ws.Range("G4").Value = 90
$ws.Range("H4").Value = 8

# Re-create the hyperlinks in display order: Parth, Piyush, Shaweta so the
# relationship ids come out as rId1/rId2/rId3 respectively.
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:parth@gmail.com", "", "", "parth@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:piyush@gmail.com", "", "", "piyush@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:shaweta@gmail.com", "", "", "shaweta@gmail.com")

# Hyperlinks.Add recolors the cell with the builtin "Hyperlink" style and
# resets its text -- reapply the original style and displayed text.
$ws.Range("Z1").Copy($ws.Range("D2"))
$ws.Range("D2").Value = "parth@gmail.com"
$ws.Range("Z1").Copy($ws.Range("D3"))
$ws.Range("D3").Value = "piyush@gmail.com"
$ws.Range("Z1").Copy($ws.Range("D4"))
$ws.Range("D4").Value = "shaweta@gmail.com"

# Drop the scratch cell used to carry the hyperlink style around.
$ws.Range("Z1").Clear()

# Match the selection left behind by the edit (cursor on H2).
$ws.Range("H2").Select()
